# PM17 Tidsregistrering for Tommy - fill in time registrations for rows 24-28
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 24
$ws.Cells.Item(24, 1).Value = "Lav SD0104 hentOmsaetning"
$ws.Cells.Item(24, 2).Value = "Software Architect"
$ws.Cells.Item(24, 3).Value = "2020-03-02"
$ws.Cells.Item(24, 4).Value = 0.36458333333333331
$ws.Cells.Item(24, 5).Value = 0.38194444444444442

# Row 25
$ws.Cells.Item(25, 1).Value = "Lav DCD0104 hentOmsatning"
$ws.Cells.Item(25, 2).Value = "Software Architect"
$ws.Cells.Item(25, 3).Value = "2020-03-02"
$ws.Cells.Item(25, 4).Value = 0.38194444444444442
$ws.Cells.Item(25, 5).Value = 0.39583333333333331

# Row 26
$ws.Cells.Item(26, 1).Value = "Review SD0103 og DCD0103"
$ws.Cells.Item(26, 2).Value = "Reviewer"
$ws.Cells.Item(26, 3).Value = "2020-03-02"
$ws.Cells.Item(26, 4).Value = 0.44791666666666669
$ws.Cells.Item(26, 5).Value = 0.4548611111111111

# Row 27
$ws.Cells.Item(27, 1).Value = "Implementering af SD0101 og SD0102"
$ws.Cells.Item(27, 2).Value = "Implenter"
$ws.Cells.Item(27, 3).Value = "2020-03-02"
$ws.Cells.Item(27, 4).Value = 0.45833333333333331
$ws.Cells.Item(27, 5).Value = 0.52083333333333337

# Row 28
$ws.Cells.Item(28, 1).Value = "Lav unit tests til OC0101"
$ws.Cells.Item(28, 2).Value = "Test Analyst"
$ws.Cells.Item(28, 3).Value = "2020-03-02"
$ws.Cells.Item(28, 4).Value = 0.52083333333333337
$ws.Cells.Item(28, 5).Value = 0.60416666666666663

# Update view: selection
$ws.Range("H85").Select()
